$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R11").Value = 372
$ws.Range("S11").Value = 41.304
$ws.Range("T11").Value = 1874
$ws.Range("U11").Value = 208.074
$ws.Range("R18").Value = 371
$ws.Range("S18").Value = 32.011
$ws.Range("T18").Value = 1955
$ws.Range("U18").Value = 168.685
$ws.Range("V18").Value = 2838.847
$ws.Range("W18").Value = 244.947
$ws.Range("X18").Value = 909.524
$ws.Range("Y18").Value = 78.477
$ws.Range("R27").Value = 380
$ws.Range("S27").Value = 54.688
$ws.Range("T27").Value = 4205
$ws.Range("U27").Value = 605.171
$ws.Range("R32").Value = 899
$ws.Range("S32").Value = 23.82
$ws.Range("T32").Value = 4848
$ws.Range("U32").Value = 128.451
$ws.Range("T43").Value = 2131
$ws.Range("U43").Value = 519.0890000000001
$ws.Range("X43").Value = 1211.565
$ws.Range("Y43").Value = 295.124
$ws.Range("R45").Value = 31
$ws.Range("S45").Value = 35.392
$ws.Range("T45").Value = 187
$ws.Range("U45").Value = 213.495
$ws.Range("V45").Value = 15
$ws.Range("W45").Value = 17.125
$ws.Range("X45").Value = 120
$ws.Range("Y45").Value = 137.002
$ws.Range("R46").Value = 1098
$ws.Range("S46").Value = 102.531
$ws.Range("T46").Value = 6830
$ws.Range("U46").Value = 637.782
$ws.Range("V46").Value = 1923.631
$ws.Range("W46").Value = 179.628
$ws.Range("X46").Value = 12029.48
$ws.Range("Y46").Value = 1123.307
$ws.Range("R48").Value = 138
$ws.Range("S48").Value = 23.825
$ws.Range("T48").Value = 875
$ws.Range("U48").Value = 151.065
$ws.Range("X48").Value = 861.9349999999999
$ws.Range("Y48").Value = 148.809
$ws.Range("R57").Value = 43
$ws.Range("S57").Value = 32.415
$ws.Range("T57").Value = 418
$ws.Range("U57").Value = 315.106
$ws.Range("V57").Value = 31.04
$ws.Range("W57").Value = 23.399
$ws.Range("X57").Value = 277.359
$ws.Range("Y57").Value = 209.085
$ws.Range("R61").Value = 29
$ws.Range("S61").Value = 5.234
$ws.Range("T61").Value = 176
$ws.Range("U61").Value = 31.765
$ws.Range("R62").Value = 2620
$ws.Range("S62").Value = 40.139
$ws.Range("T62").Value = 24526
$ws.Range("U62").Value = 375.742
$ws.Range("V62").Value = 1325.674
$ws.Range("W62").Value = 20.31
$ws.Range("X62").Value = 8387.495000000001
$ws.Range("Y62").Value = 128.498
$ws.Range("R66").Value = 5300
$ws.Range("S66").Value = 63.258
$ws.Range("X66").Value = 2324.214
$ws.Range("Y66").Value = 27.741
$ws.Range("V68").Value = 103.991
$ws.Range("W68").Value = 9.977
$ws.Range("T77").Value = 4980
$ws.Range("U77").Value = 515.509
$ws.Range("R78").Value = 0
$ws.Range("S78").Value = 0
$ws.Range("T78").Value = 20
$ws.Range("U78").Value = 58.608
$ws.Range("R84").Value = 128
$ws.Range("S84").Value = 25.922
$ws.Range("T84").Value = 1426
$ws.Range("U84").Value = 288.793
$ws.Range("V84").Value = 29.198
$ws.Range("W84").Value = 5.913
$ws.Range("X84").Value = 695.728
$ws.Range("Y84").Value = 140.898
$ws.Range("R86").Value = 2615
$ws.Range("S86").Value = 43.25
$ws.Range("T86").Value = 26042
$ws.Range("U86").Value = 430.718
$ws.Range("X86").Value = 3948.68
$ws.Range("Y86").Value = 65.309
$ws.Range("T96").Value = 1147
$ws.Range("U96").Value = 608.1
$ws.Range("V96").Value = 116.907
$ws.Range("W96").Value = 61.98
$ws.Range("X96").Value = 921.504
$ws.Range("Y96").Value = 488.55
$ws.Range("T102").Value = 2212
$ws.Range("U102").Value = 812.551
$ws.Range("V102").Value = 109.118
$ws.Range("W102").Value = 40.083
$ws.Range("X102").Value = 649.838
$ws.Range("Y102").Value = 238.71
$ws.Range("R103").Value = 23
$ws.Range("S103").Value = 36.743
$ws.Range("T103").Value = 92
$ws.Range("U103").Value = 146.97
$ws.Range("X109").Value = 34.889
$ws.Range("Y109").Value = 79.018
$ws.Range("R123").Value = 711
$ws.Range("S123").Value = 41.494
$ws.Range("T123").Value = 1850
$ws.Range("U123").Value = 107.967
$ws.Range("V123").Value = 278.605
$ws.Range("W123").Value = 16.26
$ws.Range("X123").Value = 1552.654
$ws.Range("Y123").Value = 90.614
$ws.Range("T129").Value = 143
$ws.Range("U129").Value = 26.378
$ws.Range("V129").Value = 19.332
$ws.Range("W129").Value = 3.566
$ws.Range("X129").Value = 106.833
$ws.Range("Y129").Value = 19.706
$ws.Range("T138").Value = 16727
$ws.Range("U138").Value = 441.968
$ws.Range("R139").Value = 558
$ws.Range("S139").Value = 54.724
$ws.Range("T139").Value = 3770
$ws.Range("U139").Value = 369.727
$ws.Range("X139").Value = 455.431
$ws.Range("Y139").Value = 44.665
$ws.Range("R141").Value = 1065
$ws.Range("S141").Value = 55.36
$ws.Range("T156").Value = 2980
$ws.Range("U156").Value = 545.823
$ws.Range("R157").Value = 202
$ws.Range("S157").Value = 97.16500000000001
$ws.Range("T157").Value = 1192
$ws.Range("U157").Value = 573.371
$ws.Range("V157").Value = 113.892
$ws.Range("W157").Value = 54.784
$ws.Range("X157").Value = 743.294
$ws.Range("Y157").Value = 357.536
$ws.Range("R163").Value = 2318
$ws.Range("S163").Value = 49.578
$ws.Range("T163").Value = 14951
$ws.Range("U163").Value = 319.775
$ws.Range("V163").Value = 94.631
$ws.Range("W163").Value = 2.024
$ws.Range("X163").Value = 1832.854
$ws.Range("Y163").Value = 39.201
$ws.Range("R167").Value = 372
$ws.Range("S167").Value = 36.834
$ws.Range("T167").Value = 2844
$ws.Range("U167").Value = 281.605
$ws.Range("V167").Value = 172.761
$ws.Range("W167").Value = 17.106
$ws.Range("R182").Value = 3626
$ws.Range("S182").Value = 53.413
$ws.Range("T182").Value = 36797
$ws.Range("U182").Value = 542.0410000000001
$ws.Range("X182").Value = 16617
$ws.Range("Y182").Value = 244.778
$ws.Range("R183").Value = 23891
$ws.Range("S183").Value = 72.178
$ws.Range("T183").Value = 128947
$ws.Range("U183").Value = 389.565
$ws.Range("V183").Value = 1273
$ws.Range("W183").Value = 3.846
$ws.Range("X183").Value = 33902
$ws.Range("Y183").Value = 102.422
